$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.033.85"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "1.881.98"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.79"
$ws.Range("E5").Value = "  -2.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4958"
$ws.Range("E7").Value = "  -0.61%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.36"
$ws.Range("E8").Value = "  -2.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2918"
$ws.Range("E9").Value = "  +2.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06613"
$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").Value = "1.879.86"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.89"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07198"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6659"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.90"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.839"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "30.019.93"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007782"
$ws.Range("E18").Value = "  +3.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9988"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.81"
$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("D21").Value = "2.120.99"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9983"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.762"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.596"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.148"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.19"
$ws.Range("E26").Value = "  +3.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.81"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.905"
$ws.Range("E29").Value = "  -2.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.381"
$ws.Range("E30").Value = "  -0.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.171"
$ws.Range("E31").Value = "  -1.63%  "

$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.950"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04986"
$ws.Range("E34").Value = "  -1.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.104"
$ws.Range("E35").Value = "  -2.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7014"
$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.690"
$ws.Range("E38").Value = "  -2.08%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.193"
$ws.Range("E39").Value = "  -6.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9336"
$ws.Range("E40").Value = "  -2.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01641"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.961"
$ws.Range("E42").Value = "  -2.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9992"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4181"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.33"
$ws.Range("E45").Value = "  -2.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.529"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1260"
$ws.Range("E47").Value = "  +0.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05719"
$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.33"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.194"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("E51").Value = "  +1.74%  "
